$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ncbitax")
$c = $ws.Range("I15")
$c.Interior.Color = 255
$c.Interior.Pattern = -4142
$c.Interior.ColorIndex = -4105
